$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: C1 "LF" -> "V"
$ws.Range("C1").Value = "V"

# A3: "LF Lag" -> "V Lag"
$ws.Range("A3").Value = "V Lag"

# Numeric-looking text values (B2, C2, B3, C3) need to stay as TEXT (shared
# string), matching the source workbook where these are all stored as
# strings, not numbers. Force text formatting before assigning so the
# engine doesn't coerce the numeric-looking strings into numbers, then
# clear the formatting back off so the cells keep their original (default)
# style once the text type has been locked in.
$numRange = $ws.Range("B2:C3")
$numRange.NumberFormat = "@"

$ws.Range("B2").Value = "-0.04"
$ws.Range("C2").Value = "5.81"
$ws.Range("B3").Value = "0.0"
$ws.Range("C3").Value = "-1.11***"

$numRange.ClearFormats()
